# Updates the weekly "Ají" (Hortaliza) price records for the Mapocho
# Venta Directa de Santiago market. The underlying data rows (2-18) have
# been refreshed with this week's figures; most rows keep the same shape
# but the Fecha / Variedad / Calidad / Volumen / Precio* / Unidad /
# Precio $/Kg / Kg o Unidades values move around between records.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsData = @{
    2  = @{ D = 44553; H = "Inferno";        I = "Primera"; J = 35; K = 45000;  L = 45000;  M = 45000;  N = "`$/caja 25 kilos"; P = 1800; Q = 25 }
    3  = @{ D = 44425; H = "Americana (o)";  I = "Primera"; J = 15; K = 75000;  L = 75000;  M = 75000;  N = "`$/caja 25 kilos"; P = 3000; Q = 25 }
    4  = @{ D = 44340; H = "Americana (o)";  I = "Primera"; J = 15; K = 35000;  L = 35000;  M = 35000;  N = "`$/caja 25 kilos"; P = 1400; Q = 25 }
    5  = @{ D = 44449; H = "Americana (o)";  I = "Primera"; J = 25; K = 80000;  L = 80000;  M = 80000;  N = "`$/caja 25 kilos"; P = 3200; Q = 25 }
    6  = @{ D = 44449; H = "Americana (o)";  I = "Segunda"; J = 20; K = 75000;  L = 75000;  M = 75000;  N = "`$/caja 15 kilos"; P = 5000; Q = 15 }
    7  = @{ D = 44544; H = "Inferno";        I = "Primera"; J = 12; K = 35000;  L = 35000;  M = 35000;  N = "`$/caja 25 kilos"; P = 1400; Q = 25 }
    8  = @{ D = 44319; H = "Americana (o)";  I = "Primera"; J = 20; K = 30000;  L = 30000;  M = 30000;  N = "`$/caja 25 kilos"; P = 1200; Q = 25 }
    9  = @{ D = 44474; H = "Americana (o)";  I = "Primera"; J = 18; K = 100000; L = 100000; M = 100000; N = "`$/caja 25 kilos"; P = 4000; Q = 25 }
    10 = @{ D = 44343; H = "Americana (o)";  I = "Primera"; J = 20; K = 36000;  L = 36000;  M = 36000;  N = "`$/caja 25 kilos"; P = 1440; Q = 25 }
    11 = @{ D = 44460; H = "Americana (o)";  I = "Primera"; J = 30; K = 95000;  L = 95000;  M = 95000;  N = "`$/caja 25 kilos"; P = 3800; Q = 25 }
    12 = @{ D = 44326; H = "Americana (o)";  I = "Primera"; J = 15; K = 30000;  L = 30000;  M = 30000;  N = "`$/caja 25 kilos"; P = 1200; Q = 25 }
    13 = @{ D = 44581; H = "Americana (o)";  I = "Segunda"; J = 30; K = 17000;  L = 17000;  M = 17000;  N = "`$/caja 25 kilos"; P = 680;  Q = 25 }
    14 = @{ D = 44446; H = "Americana (o)";  I = "Primera"; J = 5;  K = 78000;  L = 78000;  M = 78000;  N = "`$/caja 25 kilos"; P = 3120; Q = 25 }
    15 = @{ D = 44446; H = "Inferno";        I = "Primera"; J = 4;  K = 80000;  L = 80000;  M = 80000;  N = "`$/caja 15 kilos"; P = 5333; Q = 15 }
    16 = @{ D = 44421; H = "Americana (o)";  I = "Primera"; J = 15; K = 75000;  L = 75000;  M = 75000;  N = "`$/caja 25 kilos"; P = 3000; Q = 25 }
    17 = @{ D = 44193; H = "Americana (o)";  I = "Primera"; J = 15; K = 46000;  L = 46000;  M = 46000;  N = "`$/caja 15 kilos"; P = 3067; Q = 15 }
    18 = @{ D = 44221; H = "Americana (o)";  I = "Primera"; J = 22; K = 24000;  L = 25000;  M = 24545;  N = "`$/caja 25 kilos"; P = 982;  Q = 25 }
}

foreach ($r in $rowsData.Keys) {
    $row = $rowsData[$r]
    $ws.Cells.Item($r, 4).Value  = $row.D   # Fecha
    $ws.Cells.Item($r, 8).Value  = $row.H   # Variedad
    $ws.Cells.Item($r, 9).Value  = $row.I   # Calidad
    $ws.Cells.Item($r, 10).Value = $row.J   # Volumen
    $ws.Cells.Item($r, 11).Value = $row.K   # Precio minimo
    $ws.Cells.Item($r, 12).Value = $row.L   # Precio maximo
    $ws.Cells.Item($r, 13).Value = $row.M   # Precio promedio ponderado
    $ws.Cells.Item($r, 14).Value = $row.N   # Unidad de comercializacion
    $ws.Cells.Item($r, 16).Value = $row.P   # Precio $/Kg
    $ws.Cells.Item($r, 17).Value = $row.Q   # Kg o Unidades
}
